# LoginData.xlsx update
# - "Login" sheet (sheet1): row 2 now holds raviuser@yopmail.com / 12345678
#   instead of amrendrasadmin@yopmail.com / pass1234.
# - "LinksLogin" sheet (sheet2): row 2 keeps raviuser@yopmail.com but the
#   cell gets a plain-text number format; page setup (A4/portrait) is added.
# - Active sheet / selection state flips from LinksLogin -> Login.

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsLinks = $wb.Worksheets.Item("LinksLogin")

# ---- LinksLogin (sheet2) --------------------------------------------------
# Select it first (while it is still the active tab) so the leftover
# selection anchor matches the target state.
$wsLinks.Activate()
$wsLinks.Range("B14").Select()

# A2 keeps the same hyperlink text, but is formatted as text (numFmt 49)
$wsLinks.Range("A2").NumberFormat = "@"

$wsLinks.PageSetup.PaperSize = 9
$wsLinks.PageSetup.Orientation = 1

# ---- Login (sheet1) --------------------------------------------------------
$wsLogin.Activate()

$wsLogin.Range("A2").Value = "raviuser@yopmail.com"
$wsLogin.Range("A2").Borders.LineStyle = 1

$wsLogin.Range("B2").Value = 12345678

$wsLogin.Range("A21").Select()
